# LOB1049.docx edit: rotates the course-description paragraph contents.
#
# The document's section values are cyclically shifted one slot "down":
#   Objetivos value             <- old Programa-resumido value
#   Docente(s) list item        <- old Objetivos value
#   Programa-resumido value     <- old Programa value
#   Programa value              <- old Metodo value
#   Metodo value                <- old Criterio value
#   Criterio value              <- old Norma-de-recuperacao value
#   Norma-de-recuperacao value  <- old Bibliografia value
#   Bibliografia value          <- old Docente(s) list item
#
# Because every destination value equals some *other* paragraph's original
# value, writing the new values directly (old-text -> new-text, using the
# literal captured strings) would create duplicate text in the document
# partway through, and later Find operations could then latch onto the
# wrong (already-rewritten) location. To avoid that, we do this in two
# clean passes:
#   Pass 1: replace each of the 8 original values with a unique synthetic
#           placeholder token (each original value occurs exactly once, so
#           this pass is unambiguous).
#   Pass 2: replace each placeholder token with the real destination text
#           (the placeholders are unique and can't collide with anything),
#           so ordering no longer matters.

$d = $word.ActiveDocument
$vt = [char]11

function Replace-Text($oldText, $newText) {
    $range = $d.Content
    $ok = $range.Find.Execute(
        $oldText,   # FindText
        $true,      # MatchCase
        $false,     # MatchWholeWord
        $false,     # MatchWildcards
        $false,     # MatchSoundsLike
        $false,     # MatchAllWordForms
        $true,      # Forward
        1,          # Wrap = wdFindContinue
        $false,     # Format
        $newText,   # ReplaceWith
        2           # Replace = wdReplaceOne
    )
    if (-not $ok) {
        throw "Find/Replace failed for: $oldText"
    }
}

# ---- capture the 8 original values (pre-edit state) ----

$objetivos_old = "Fornecer aos alunos os conceitos básicos de Estatística Multivariada assim como sua aplicação nos estudos de fenômenos onde vários componentes se comportam de forma correlacionada."

$docente_old = "4894221 - Mariana Pereira de Melo"

$resumido_old = "Probabilidade: Vetor de variáveis aleatórias, Distribuição conjunta/marginal, Esperança e variância condicional/marginal. Estatística: Regressão Logística simples, Teste Qui-Quadrado, Testes de normalidade, Testes não-paramétricos. Técnicas Multivariadas: Gráficos multivariados, Regressão Linear Múltipla, Regressão Logística Múltipla, Análise de Variância Múltipla; Análise de agrupamento; Análise de componentes principais; Análise fatorial; Análise discriminante e Análise de correspondência"

$programa_old = "Probabilidade: Vetor de variáveis aleatórias, Distribuição conjunta/marginal, Esperança e Variância condicional/marginal." + $vt + "Estatística: Regressão Logística simples (coeficiente de associação, sensitividade e especificidade, risco relativo, razão de chances), Teste Qui-Quadrado (testes de aderência, homogeneidade e independência), Teste de normalidade (Shapiro-Wilk, Teste de Kolmogorov-Smirnov), Testes não-paramétricos para amostras pareadas e independentes." + $vt + "Técnicas Multivariadas: Gráficos multivariados, Regressão Linear Múltipla, Regressão Logística Múltipla, Análise de Variância Múltipla; Análise de agrupamento; Análise de componentes principais; Análise fatorial; Análise discriminante e Análise de correspondência"

$metodo_old = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: N = (N1+...+Nn)/n"

$criterio_old = "NF≥ 5,0."

$norma_old = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."

$biblio_old = "G.C. Runger, D. Montgomery. Estatística aplicada e probabilidade para engenheiros. São Paulo: Ed. LTC, 2009. " + $vt + $vt + "D. C. Montgomery, E. A. Peck, G. G. Vining, Introduction to Linear Regression Analysis, 4th ed., Hoboken: John Wiley, 2006." + $vt + $vt + "W. J. Conover, Practical Nonparametric Statistics, 3rd ed., New York: John Wiley d Sons, 1999." + $vt + $vt + "R. A. Johnson, D. W. Wichern, Applied Multivariate Statistical Analysis, 6th ed., New Jersey: Prentice Hall, 2007."

# ---- Pass 1: stamp each original location with a unique placeholder ----

Replace-Text $objetivos_old "@@SLOT_OBJETIVOS@@"
Replace-Text $docente_old    "@@SLOT_DOCENTE@@"
Replace-Text $resumido_old   "@@SLOT_RESUMIDO@@"
Replace-Text $programa_old   "@@SLOT_PROGRAMA@@"
Replace-Text $metodo_old     "@@SLOT_METODO@@"
Replace-Text $criterio_old   "@@SLOT_CRITERIO@@"
Replace-Text $norma_old      "@@SLOT_NORMA@@"
Replace-Text $biblio_old     "@@SLOT_BIBLIO@@"

# ---- Pass 2: fill each placeholder with its final destination content ----

Replace-Text "@@SLOT_OBJETIVOS@@" $resumido_old
Replace-Text "@@SLOT_DOCENTE@@"   $objetivos_old
Replace-Text "@@SLOT_RESUMIDO@@"  $programa_old
Replace-Text "@@SLOT_PROGRAMA@@"  $metodo_old
Replace-Text "@@SLOT_METODO@@"    $criterio_old
Replace-Text "@@SLOT_CRITERIO@@"  $norma_old
Replace-Text "@@SLOT_NORMA@@"     $biblio_old
Replace-Text "@@SLOT_BIBLIO@@"    $docente_old

Write-Output "LOB1049 rotation applied"
